$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-152: refresh submission timestamp (column A) and formatted date (column K)
# to reflect the latest resubmission of the existing survey responses.
for ($r = 2; $r -le 152; $r++) {
    $ws.Cells.Item($r, 1).Value = 1527859196
    $ws.Cells.Item($r, 11).Value = "Fri_Jun__1_09:19:56_EDT_2018"
}

# Rows 153-177: newly reported TH Real Estate locations (ParentCorpID groups 8 & 9)
# plus refreshed "All locations" summary rows for every ParentCorpID group, and the
# latest single-location resubmission row.
$newRows = @(
    @(1527859196, 8, "Kerry_Ireland", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "Grasse_France", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "Guanajuato_Mexico", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "LapulapuCity_Cebu_Philipines", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "Singapore", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "Lake_Zuirich_Illinois_USA", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "Manaus_Brazil", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "Western_Cape_South_Africa", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "Sydney_Australia", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "Bangalore_india", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "Los_Angeles_USA", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "Pontevedra_Spain", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 9, "801_Brickell_Miami_USA", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 9, "Xanadu_Centre_Madrid_Spain", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 9, "Ginza_1_Chome_Tokyo_Japan", "-", "250", "100", "90", "80", "0.034", "Clean Room Manufacturing, R&D", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 1, "All locations", "Manufacturing", "-", "-", "-", "-", "-", "-", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 2, "All locations", "Manufacturing", "-", "-", "-", "-", "-", "-", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 3, "All locations", "Manufacturing", "-", "-", "-", "-", "-", "-", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 4, "All locations", "Manufacturing", "-", "-", "-", "-", "-", "-", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 5, "All locations", "Manufacturing", "-", "-", "-", "-", "-", "-", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 6, "All locations", "Manufacturing", "-", "-", "-", "-", "-", "-", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 7, "All locations", "Manufacturing", "-", "-", "-", "-", "-", "-", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 8, "All locations", "Manufacturing", "-", "-", "-", "-", "-", "-", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527859196, 9, "All locations", "Manufacturing", "-", "-", "-", "-", "-", "-", "Fri_Jun__1_09:19:56_EDT_2018"),
    @(1527860379, 9, " All locations ", " Beverage ", " - ", " - ", " - ", " - ", " - ", " - ", " Fri_Jun__1_09:39:39_EDT_2018")
)

$startRow = 153
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowData = $newRows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $val = $rowData[$c]
        $colNum = $c + 1
        $cell = $ws.Cells.Item($r, $colNum)
        # Numeric-looking strings (e.g. "250", "0.034") must be forced to Text format,
        # otherwise Excel auto-converts them to real numbers on assignment.
        if ($val -is [string] -and $val -match "^\s*-?[0-9]*\.?[0-9]+\s*$") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
    }
}
